$d = $word.ActiveDocument

# Helper-ish pattern used throughout: locate a unique phrase with
# Find, then Expand(4) [wdParagraph] the collapsed range to obtain the
# full enclosing paragraph (including its trailing paragraph mark), so
# the whole paragraph's runs can be merged/replaced as one block.

# ------------------------------------------------------------------
# Edit 1: merge the three runs that make up the "Halton Catholic..."
# sentence into a single run (text content is unchanged).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("it relies on the Halton Catholic District School Board", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$full = $d.Range($para.Start, $para.End - 1)
$full.Text = "One of the challenges in allocating EAs to schools is being able to weigh relevant criteria during the decision-making process. Furthermore, the method for making these decisions must be replicable and transparent. Two criteria that play a major role in determining the allotment of EAs, are the needs of the student and the finite number of available EAs. Weighing these factors " + [char]0x201C + "by hand" + [char]0x201D + " or by using only one" + [char]0x2019 + "s intuition, can be an opaque and onerous decision-making approach. Therefore, the goal of the EA algorithm is to balance these two factors using computational techniques that provide an initial decision for the allotment of EAs. Again, this is only the first pass at the decision-making process. After the algorithm provides results, the special education staff examine the numbers and adjust according to their expertise. The algorithm can only work if there is a way to quantify student need. For this, it relies on the Halton Catholic District School Board" + [char]0x2019 + "s independence rubric."

# ------------------------------------------------------------------
# Edit 2: merge the nine runs that make up the "independence rubric
# assesses..." paragraph into a single run (text content unchanged).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("independence rubric assesses special-needs students", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$full = $d.Range($para.Start, $para.End - 1)
$full.Text = "The independence rubric assesses special-needs students in seven domains: health/medical, safety, adaptive functioning, communication, social/emotional, academic, and community/leisure/work. These domains are further broken down into more specific subdomains. For example, under adaptive functioning, the subdomains include toileting, feeding, dressing, mobility, and personal hygiene. For each subdomain students are categorized as either level 1, 2, 3, or 4. The four levels represent how much support a student requires in any given area. For example, Level 1 indicates that much support is required, whereas level 4 indicates that no support is required. The categorizations provide a foundation for the EA algorithm in terms of the quantification of student need."

# ------------------------------------------------------------------
# Edit 3: merge the seven runs making up the "This mapping of
# categorical values..." paragraph into a single run (text unchanged,
# trailing space retained).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("This mapping of categorical values", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$full = $d.Range($para.Start, $para.End - 1)
$full.Text = "This mapping of categorical values to numbers in the unit interval is not arbitrary. One can imagine that a student in a particular domain who is a level 1 (on average), may require 80% of an EA" + [char]0x2019 + "s time (or .8 of an EA). Therefore, the mappings have a quantifiable meaning and these particular mappings are chosen based on consultation with special education staff. As described in the following section, these particular values may get adjusted up or down by the algorithm and therefore what is most important when choosing a mapping is the relative space between values. That is, one must consider how much more support a level 1 should generate compared to a level 2, and so on. Following the mappings, the algorithm begins an iterative process of balancing student need with the finite number of available EAs. "

# ------------------------------------------------------------------
# Edit 4: the "An estimate of the support..." paragraph is re-split
# into three runs, with new wording describing the grouping/averaging
# feature (averaging within domain, then across domains).
# ------------------------------------------------------------------
$seg1 = "An estimate of the support required for a particular student is given by "
$seg2 = "first averaging the values within each domain (yielding domain-specific estimates), and then averaging across these the domain-specific estimates (yielding a grand average). "
$seg3 = "Since the mapped values represent proportions, for any given school, the sum of the mean values across students represents the number of EAs required for that school. In this manner, the algorithm initially gives the full amount of support to students (as determined by the mapped values) with no restriction imposed by the finite number of available EAs. In the event of over allocation, the algorithm returns to the mapping stage and decrements the mapped values by a small amount (i.e., .001). For example, consider the following downward adjustments:"

$rng = $d.Content
$rng.Find.Execute("An estimate of the support required", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$pStart = $para.Start

$full = $d.Range($pStart, $para.End - 1)
$full.Text = $seg1 + $seg2 + $seg3

# Re-split the single merged run into three runs matching the target
# paragraph structure (same formatting throughout).
$r1 = $d.Range($pStart, $pStart + $seg1.Length)
$r1.Text = $seg1
$r2 = $d.Range($pStart + $seg1.Length, $pStart + $seg1.Length + $seg2.Length)
$r2.Text = $seg2
$r3 = $d.Range($pStart + $seg1.Length + $seg2.Length, $pStart + $seg1.Length + $seg2.Length + $seg3.Length)
$r3.Text = $seg3

# ------------------------------------------------------------------
# Edit 5: merge the three runs after the leading <w:tab/> in the
# "The algorithm finishes within seconds..." paragraph into one run.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The algorithm finishes within seconds", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$textStart = $rng.Start
$para = $d.Range($rng.Start, $rng.Start)
$para.Expand(4) | Out-Null
$full = $d.Range($textStart, $para.End - 1)
$full.Text = "The algorithm finishes within seconds and produces a number of reports for the special education staff. The reports detail how much support each student was assigned (in terms of a proportion of an EAs time) and how many EAs are allotted to each school. The algorithm should be rerun using different mapping parameters until the results produce a good fit with selected cases as determined by special education experts. "
